$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 32.99293476771223),
    @(3, 39.81795260052328),
    @(4, 32.8466741743693),
    @(5, 30.89634255962685),
    @(6, 30.76680808256089),
    @(7, 26.39873364736599),
    @(8, 24.86969599478817),
    @(9, 26.38991189248756),
    @(10, 23.94022531284983),
    @(11, 27.25960600151734),
    @(12, 24.75452262489402),
    @(13, 20.46655330024399),
    @(14, 21.05394258799927),
    @(15, 22.20226048738496),
    @(16, 19.46990463595043),
    @(17, 17.26670102626844),
    @(18, 21.60091681413137),
    @(19, 15.00795634773498),
    @(20, 20.09772569057796),
    @(21, 19.87328033220207),
    @(22, 17.24699887520443),
    @(23, 18.39964834137456),
    @(24, 16.9082392868022),
    @(25, 16.73228708559159),
    @(26, 13.76641543155043),
    @(27, 17.51973174560494),
    @(28, 18.69396516350307),
    @(29, 14.47722281418032),
    @(30, 13.18337501373389),
    @(31, 10.57188800643623),
    @(32, 15.11393900861776),
    @(33, 18.11137515097056),
    @(34, 17.68118765067433),
    @(35, 17.31761244851236),
    @(36, 14.39128494512975),
    @(37, 13.73328375738535),
    @(38, 14.79663574215394),
    @(39, 15.40097082537302),
    @(40, 12.75473845065363),
    @(41, 12.72654919580208),
    @(42, 15.51420423065211),
    @(43, 15.17978664506111),
    @(44, 15.40841250596976),
    @(45, 15.90846253914941),
    @(46, 16.43305033557174),
    @(47, 15.6487362404888),
    @(48, 13.85595256237559),
    @(49, 12.96146455252259),
    @(50, 15.12808471808341),
    @(51, 11.78453101333216),
    @(52, 13.48623860161457),
    @(53, 13.34904515479656),
    @(54, 11.10340493431545),
    @(55, 11.92007380866352),
    @(56, 13.23724683868024),
    @(57, 13.67155087160168),
    @(58, 13.10482456713882),
    @(59, 11.22527658932455),
    @(60, 10.82300621784216),
    @(61, 10.55655299881593),
    @(62, 11.93905418264083),
    @(63, 10.61118378234043),
    @(64, 9.381881871243422),
    @(65, 12.4443902414383),
    @(66, 10.79018172918848),
    @(67, 12.2679289839925),
    @(68, 9.873804600184599),
    @(69, 11.33949641835103),
    @(70, 12.26926047755127),
    @(71, 13.33784879145171),
    @(72, 10.33215490030048),
    @(73, 12.671808382455),
    @(74, 11.19262517769556),
    @(75, 15.48787191038721),
    @(76, 13.38839581231926),
    @(77, 12.90626170728282),
    @(78, 12.5276258280859),
    @(79, 13.19894011309694),
    @(80, 14.5247248952989),
    @(81, 12.93778939155904),
    @(82, 13.88548470606194),
    @(83, 12.77950087836109),
    @(84, 13.32780306795826),
    @(85, 15.91066099610471),
    @(86, 13.44500137144931)
)

foreach ($pair in $values) {
    $row = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($row, 1).Value = $val
}
